# Update "想去人数" (want-to-go count) figures in column F for the
# 展览 (Exhibitions) sheet and the aggregated 全部类型 (All types) sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 605
$wsExhibit.Range("F5").Value = 363
$wsExhibit.Range("F6").Value = 1906
$wsExhibit.Range("F7").Value = 103

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 605
$wsAll.Range("F5").Value = 363
$wsAll.Range("F10").Value = 1906
$wsAll.Range("F11").Value = 103
